# Insert a new weekly price record for "Acelga" (Terminal Hortofrutícola Agro
# Chillán) as row 154, pushing the existing rows 154-183 down to 155-184.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("154:154").Insert()

$ws.Cells.Item(154, 1).Value = 7
$ws.Cells.Item(154, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(154, 3).Value = "Ñuble"
$ws.Cells.Item(154, 4).Value = 44522
$ws.Cells.Item(154, 5).Value = 16
$ws.Cells.Item(154, 6).Value = 100112009
$ws.Cells.Item(154, 7).Value = "Acelga"
$ws.Cells.Item(154, 8).Value = "Sin especificar"
$ws.Cells.Item(154, 9).Value = "Primera"
$ws.Cells.Item(154, 10).Value = 100
$ws.Cells.Item(154, 11).Value = 350
$ws.Cells.Item(154, 12).Value = 400
$ws.Cells.Item(154, 13).Value = 375
$ws.Cells.Item(154, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(154, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(154, 16).Value = 375
$ws.Cells.Item(154, 17).Value = 1
$ws.Cells.Item(154, 18).Value = "Hortaliza"
